# Update the "ランサーズ" sheet: refresh fetch timestamps to 2026-01-26 12:43:00
# and append newly scraped listings, shifting existing rows down as needed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Keep a reference to an existing hyperlink-styled cell so new hyperlink cells
# can reuse the same cell style (avoids creating a duplicate style entry).
$styleSource = $ws.Range("F2")

# Drop all existing hyperlinks up front; they will be re-created below, in final
# row order, once every cell holds its target value.
$ws.Hyperlinks.Delete()

# Write the full final table for rows 2..9 (existing rows refreshed in place,
# new rows appended/inserted in their final position).
# Row 2
$ws.Cells.Item(2,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(2,2).Value = "製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5473648"
$ws.Cells.Item(2,7).Value = 390
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Cells.Item(3,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(3,2).Value = "【急募】対話型AI WebアプリMVP開発エンジニア募集"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5478844"
$ws.Cells.Item(3,7).Value = 378
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発 ◇アプリ"

# Row 4
$ws.Cells.Item(4,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(4,2).Value = "【日々の売上の集計自動化するツール】csvDLL/スプレッドシートに入力/売上管理シート仕様変更"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5479251"
$ws.Cells.Item(4,7).Value = 183
$ws.Cells.Item(4,8).Value = "◆ツール,自動化 ◇管理"

# Row 5
$ws.Cells.Item(5,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(5,2).Value = "自動化システム"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Cells.Item(5,7).Value = 110
$ws.Cells.Item(5,8).Value = "◆自動化"

# Row 6
$ws.Cells.Item(6,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(6,2).Value = "【業務委託/パートタイム】 フルスタックエンジニア(ノーコード/JavaScript)英語案件"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5479148"
$ws.Cells.Item(6,7).Value = 78
$ws.Cells.Item(6,8).Value = "★Java"

# Row 7
$ws.Cells.Item(7,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(7,2).Value = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5468432"
$ws.Cells.Item(7,7).Value = 75
$ws.Cells.Item(7,8).Value = "◆開発"

# Row 8
$ws.Cells.Item(8,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(8,2).Value = "【急募】Webアプリ超簡易実装|5,000円でお手伝いします!"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5479368"
$ws.Cells.Item(8,7).Value = 30
$ws.Cells.Item(8,8).Value = "◇アプリ"

# Row 9
$ws.Cells.Item(9,1).Value = "2026-01-26 12:43:00"
$ws.Cells.Item(9,2).Value = "【医療保険】オンライン資格確認・請求端末セットアップ依頼"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5478715"
$ws.Cells.Item(9,7).Value = 13

# Row 9 ("medical insurance" listing) has no skill-summary column, by design.
$ws.Cells.Item(9,8).ClearContents()

# Re-create the URL hyperlinks in row order so relationship ids line up again.
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5473648")
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5478844")
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5479251")
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5477084")
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5479148")
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5468432")
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5479368")
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5478715")

# Re-apply the shared "Hyperlink" cell style to every URL cell so no duplicate
# style entries get created in styles.xml.
$styleSource.Copy()
$ws.Range("F2:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Make sure the sheet dimension / selection reflect the now-larger used range.
$ws.Range("A1").Select()
